$wb = $excel.ActiveWorkbook

# --- Colaboradores: add "Nome" header, push "Adelma" down to row 6 ---
$wsColab = $wb.Worksheets.Item("Colaboradores")
$wsColab.Range("A1").Value = "Nome"
$wsColab.Range("A6").Value = "Adelma"

# --- Turno: insert "Turno" header at the top, shifting the shift codes down ---
$wsTurno = $wb.Worksheets.Item("Turno")
$wsTurno.Range("A4").Value = "N18"
$wsTurno.Range("A3").Value = "T15"
$wsTurno.Range("A2").Value = "M6"
$wsTurno.Range("A1").Value = "Turno"

# --- Re-create the user's click-through of the workbook: select cells on a
#     few sheets, then land on ColaboradoresCargaHoraria as the final active tab. ---
$wsArea = $wb.Worksheets.Item("Area")
[void]$wsArea.Activate()
[void]$wsArea.Range("C1").Select()

[void]$wsColab.Activate()
[void]$wsColab.Range("A7").Select()

[void]$wsTurno.Activate()
[void]$wsTurno.Range("F8").Select()

$wsCarga = $wb.Worksheets.Item("ColaboradoresCargaHoraria")
[void]$wsCarga.Activate()
